# Calibrate electricity sector results against ReEDS
#
# Updates the "Share of cost effective capacity built in a single year"
# assumption for onshore wind (row 7) from 33% to 20% across all forecast
# years (2021-2046, columns B:AE) on the
# "CSC-CSCSoCECBiaSY" sheet, and leaves the workbook positioned/selected
# on that sheet & range, matching the author's last-saved view state.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")

# Onshore wind ("onshore wind es") is row 7; update the share of cost
# effective capacity built in a single year from 0.33 to 0.2 for every
# year column (B:AE).
$wsData.Range("B7:AE7").Value = 0.2

# Bring this sheet to the front and leave the same cells selected that
# the author had selected when they saved, so the view state matches.
[void]$wsData.Activate()
[void]$wsData.Range("B7:AE7").Select()
